# Update "想去人数" (interest count) figures in column F across the four
# sheets of the workbook, reflecting refreshed stats from the data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 294
$ws1.Range("F4").Value = 1223
$ws1.Range("F5").Value = 343
$ws1.Range("F6").Value = 312
$ws1.Range("F7").Value = 3630
$ws1.Range("F9").Value = 737
$ws1.Range("F10").Value = 1249
$ws1.Range("F11").Value = 324
$ws1.Range("F12").Value = 212
$ws1.Range("F13").Value = 725
$ws1.Range("F14").Value = 142
$ws1.Range("F15").Value = 149
$ws1.Range("F16").Value = 2039
$ws1.Range("F20").Value = 322
$ws1.Range("F23").Value = 265

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 76
$ws2.Range("F13").Value = 221
$ws2.Range("F23").Value = 50

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6389
$ws3.Range("F5").Value = 299

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6389
$ws4.Range("F5").Value = 299
$ws4.Range("F12").Value = 294
$ws4.Range("F13").Value = 1223
$ws4.Range("F14").Value = 343
$ws4.Range("F18").Value = 312
$ws4.Range("F19").Value = 3630
$ws4.Range("F24").Value = 76
$ws4.Range("F25").Value = 737
$ws4.Range("F26").Value = 1249
$ws4.Range("F27").Value = 324
$ws4.Range("F28").Value = 221
$ws4.Range("F29").Value = 212
$ws4.Range("F30").Value = 725
$ws4.Range("F31").Value = 142
$ws4.Range("F32").Value = 149
$ws4.Range("F34").Value = 2039
$ws4.Range("F40").Value = 322
$ws4.Range("F49").Value = 50
$ws4.Range("F50").Value = 265
